$wb = $excel.ActiveWorkbook

# --- Cases sheet ---
$ws = $wb.Worksheets.Item("Cases")
$ws.Range("AB51").Value = 26139
$ws.Range("Y51").Value = 1707
$ws.Range("AB52").Value = 26459
$ws.Range("Y52").Value = 1723
$ws.Range("AB53").Value = 26773
$ws.Range("Y53").Value = 1740
$ws.Range("AB54").Value = 27075
$ws.Range("Y54").Value = 1758
$ws.Range("AB55").Value = 27397
$ws.Range("Y55").Value = 1767
$ws.Range("AB56").Value = 27595
$ws.Range("Y56").Value = 1774
$ws.Range("AB57").Value = 27787
$ws.Range("Y57").Value = 1790
$ws.Range("AB58").Value = 27944
$ws.Range("Y58").Value = 1796
$ws.Range("AB59").Value = 28146
$ws.Range("Y59").Value = 1801
$ws.Range("AB60").Value = 28258

# --- Fatalities sheet ---
$ws = $wb.Worksheets.Item("Fatalities")
$ws.Range("AB20").Value = 12
$ws.Range("Y20").Value = 0
$ws.Range("AB21").Value = 20
$ws.Range("Y21").Value = 1
$ws.Range("AB22").Value = 28
$ws.Range("Y22").Value = 1
$ws.Range("AB23").Value = 33
$ws.Range("Y23").Value = 0
$ws.Range("AB24").Value = 46
$ws.Range("Y24").Value = 1
$ws.Range("AB25").Value = 52
$ws.Range("Y25").Value = 2
$ws.Range("AB26").Value = 72
$ws.Range("Y26").Value = 1
$ws.Range("AB27").Value = 90
$ws.Range("Y27").Value = 3
$ws.Range("AB28").Value = 105
$ws.Range("Y28").Value = 1
$ws.Range("AB29").Value = 133
$ws.Range("Y29").Value = 2
$ws.Range("AB30").Value = 152
$ws.Range("Y30").Value = 1
$ws.Range("AB31").Value = 184
$ws.Range("Y31").Value = 1
$ws.Range("AB32").Value = 233
$ws.Range("Y32").Value = 2
$ws.Range("AB33").Value = 267
$ws.Range("Y33").Value = 4
$ws.Range("AB34").Value = 311
$ws.Range("Y34").Value = 4
$ws.Range("AB35").Value = 348
$ws.Range("Y35").Value = 6
$ws.Range("AB36").Value = 408
$ws.Range("Y36").Value = 4
$ws.Range("AB37").Value = 463
$ws.Range("Y37").Value = 2
$ws.Range("AB38").Value = 522
$ws.Range("Y38").Value = 3
$ws.Range("AB39").Value = 580
$ws.Range("Y39").Value = 7
$ws.Range("AB40").Value = 636
$ws.Range("Y40").Value = 4
$ws.Range("AB41").Value = 695
$ws.Range("Y41").Value = 2
$ws.Range("AB42").Value = 746
$ws.Range("Y42").Value = 4
$ws.Range("AB43").Value = 798
$ws.Range("Y43").Value = 3
$ws.Range("AB44").Value = 856
$ws.Range("Y44").Value = 8
$ws.Range("AB45").Value = 916
$ws.Range("Y45").Value = 3
$ws.Range("AB46").Value = 962
$ws.Range("Y46").Value = 6
$ws.Range("AB47").Value = 1006
$ws.Range("Y47").Value = 6
$ws.Range("AB48").Value = 1036
$ws.Range("Y48").Value = 2
$ws.Range("AB49").Value = 1088
$ws.Range("Y49").Value = 6
$ws.Range("AB50").Value = 1111
$ws.Range("Y50").Value = 2
$ws.Range("AB51").Value = 1148
$ws.Range("Y51").Value = 1
$ws.Range("AB52").Value = 1200
$ws.Range("Y52").Value = 0
$ws.Range("AB53").Value = 1241
$ws.Range("Y53").Value = 1
$ws.Range("AB54").Value = 1282
$ws.Range("Y54").Value = 4
$ws.Range("AB55").Value = 1317
$ws.Range("Y55").Value = 5
$ws.Range("AB56").Value = 1333
$ws.Range("Y56").Value = 3
$ws.Range("AB57").Value = 1366
$ws.Range("Y57").Value = 3
$ws.Range("AB58").Value = 1399
$ws.Range("Y58").Value = 4
$ws.Range("AB59").Value = 1418
$ws.Range("Y59").Value = 3
$ws.Range("AB60").Value = 1427
$ws.Range("Y60").Value = 3

# --- Hospitalized sheet ---
$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Range("AB58").Value = 1424
$ws.Range("Y58").Value = 75
$ws.Range("AB59").Value = 1363
$ws.Range("Y59").Value = 77
$ws.Range("Y60").Value = 76

# --- ICU sheet ---
$ws = $wb.Worksheets.Item("ICU")
$ws.Range("AB53").Value = 296
$ws.Range("Y53").Value = 17
$ws.Range("AB54").Value = 290
$ws.Range("Y54").Value = 16
$ws.Range("AB55").Value = 277
$ws.Range("Y55").Value = 18
$ws.Range("AB56").Value = 269
$ws.Range("Y56").Value = 16
$ws.Range("AB57").Value = 256
$ws.Range("Y57").Value = 15
$ws.Range("AB58").Value = 244
$ws.Range("Y58").Value = 14
$ws.Range("AB59").Value = 209
$ws.Range("Y59").Value = 13
$ws.Range("AB60").Value = 203
$ws.Range("Y60").Value = 13

# --- Ventilated sheet ---
$ws = $wb.Worksheets.Item("Ventilated")
$ws.Range("Y60").Value = 11

# --- Released sheet ---
$ws = $wb.Worksheets.Item("Released")
$ws.Range("AB60").Value = 4255
$ws.Range("Y60").Value = 211

Write-Host "Applied all changes"